$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3237.5
$ws.Range("I17").Value = 950
$ws.Range("K17").Value = 2850
$ws.Range("M17").Value = -2682
$ws.Range("H86").Value = 5666.6665
$ws.Range("J86").Value = 5666.6665
$ws.Range("L86").Value = 5666.6665
$ws.Range("N86").Value = -7912.6665
$ws.Range("H89").Value = 5666.6665
$ws.Range("J89").Value = 5666.6665
$ws.Range("L89").Value = 28333.3325
$ws.Range("N89").Value = -39565.3325
$ws.Range("H100").Value = 3327.6
$ws.Range("J100").Value = 2995
$ws.Range("L100").Value = 2995
$ws.Range("N100").Value = -4077
$ws.Range("H138").Value = 6792.222
$ws.Range("J138").Value = 7410
$ws.Range("L138").Value = 22230
$ws.Range("N138").Value = -32510

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1644
$ws.Range("I2").Value = 1644
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1644
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1531
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0
$ws.Range("H26").Value = 2500
$ws.Range("I26").Value = 2500
$ws.Range("K26").Value = 2500
$ws.Range("M26").Value = -2170
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("N27").Value = 0
$ws.Range("H32").Value = 16050
$ws.Range("I32").Value = 13344.444
$ws.Range("K32").Value = 13344.444
$ws.Range("M32").Value = -13057.444
$ws.Range("H102").Value = 2906.0476
$ws.Range("I102").Value = 2484.2727
$ws.Range("J102").Value = 3370
$ws.Range("K102").Value = 2484.2727
$ws.Range("L102").Value = 3370
$ws.Range("M102").Value = -862.2727
$ws.Range("N102").Value = -6614
$ws.Range("H116").Value = 1644
$ws.Range("I116").Value = 1644
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1644
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 650
$ws.Range("H122").Value = 2779.5625
$ws.Range("I122").Value = 2746.1428
$ws.Range("J122").Value = 3013.5
$ws.Range("K122").Value = 8238.428400000001
$ws.Range("L122").Value = 9040.5
$ws.Range("M122").Value = -5788.428400000001
$ws.Range("N122").Value = -13940.5
$ws.Range("H132").Value = 2526.7
$ws.Range("I132").Value = 1295.1666
$ws.Range("J132").Value = 4374
$ws.Range("K132").Value = 3885.4998
$ws.Range("L132").Value = 13122
$ws.Range("M132").Value = -1355.4998
$ws.Range("N132").Value = -18182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1644
$ws.Range("I3").Value = 1644
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1644
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1530
$ws.Range("H25").Value = 17436.334
$ws.Range("I25").Value = 1154.5
$ws.Range("K25").Value = 1154.5
$ws.Range("M25").Value = -919.5
$ws.Range("H86").Value = 9083.166999999999
$ws.Range("I86").Value = 3999.6667
$ws.Range("K86").Value = 3999.6667
$ws.Range("M86").Value = -2876.6667
$ws.Range("H89").Value = 9083.166999999999
$ws.Range("I89").Value = 3999.6667
$ws.Range("K89").Value = 19998.3335
$ws.Range("M89").Value = -14382.3335
$ws.Range("H107").Value = 2095
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H134").Value = 1296.091
$ws.Range("I134").Value = 1361.8889
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 4085.6667
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -1550.6667
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2232.4666
$ws.Range("I31").Value = 2037.7273
$ws.Range("J31").Value = 2768
$ws.Range("K31").Value = 2037.7273
$ws.Range("L31").Value = 2768
$ws.Range("M31").Value = -1742.7273
$ws.Range("N31").Value = -3358
$ws.Range("H34").Value = 2232.4666
$ws.Range("I34").Value = 2037.7273
$ws.Range("J34").Value = 2768
$ws.Range("K34").Value = 2037.7273
$ws.Range("L34").Value = 2768
$ws.Range("M34").Value = -1835.7273
$ws.Range("N34").Value = -3172
$ws.Range("H41").Value = 23666.334
$ws.Range("J41").Value = 23666.334
$ws.Range("L41").Value = 23666.334
$ws.Range("N41").Value = -24522.334
$ws.Range("H132").Value = 4014.6
$ws.Range("J132").Value = 4014.6
$ws.Range("L132").Value = 12043.8
$ws.Range("N132").Value = -17103.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1300
$ws.Range("I132").Value = 1300
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11700
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -9170
$ws.Range("H141").Value = 6719.5
$ws.Range("I141").Value = 6966.222
$ws.Range("J141").Value = 4499
$ws.Range("K141").Value = 20898.666
$ws.Range("L141").Value = 13497
$ws.Range("M141").Value = -15718.666
$ws.Range("N141").Value = -23857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2201.9
$ws.Range("I126").Value = 2189.5715
$ws.Range("J126").Value = 2230.6667
$ws.Range("K126").Value = 6568.7145
$ws.Range("L126").Value = 6692.000100000001
$ws.Range("M126").Value = -4098.7145
$ws.Range("N126").Value = -11632.0001
$ws.Range("H132").Value = 3082.5
$ws.Range("I132").Value = 2124.25
$ws.Range("K132").Value = 6372.75
$ws.Range("M132").Value = -3842.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2781.4707
$ws.Range("I68").Value = 2767.8125
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 2767.8125
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -2018.8125
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 2781.4707
$ws.Range("I71").Value = 2767.8125
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 13839.0625
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -10095.0625
$ws.Range("N71").Value = -22488
$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676
$ws.Range("H132").Value = 3541.7334
$ws.Range("I132").Value = 3309.4546
$ws.Range("J132").Value = 4180.5
$ws.Range("K132").Value = 9928.363799999999
$ws.Range("L132").Value = 12541.5
$ws.Range("M132").Value = -7398.363799999999
$ws.Range("N132").Value = -17601.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 34994.375
$ws.Range("J54").Value = 34994.375
$ws.Range("L54").Value = 34994.375
$ws.Range("N54").Value = -36034.375
$ws.Range("H112").Value = 100000
$ws.Range("J112").Value = 100000
$ws.Range("L112").Value = 100000
$ws.Range("N112").Value = -102954
$ws.Range("H132").Value = 3002.3044
$ws.Range("I132").Value = 1475.4445
$ws.Range("K132").Value = 4426.333500000001
$ws.Range("M132").Value = -1896.333500000001
